$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Defaults" section: the Email Domain used to build the calculated email
# addresses moves from the old gmail.com account over to the new
# speridian.com one.
$ws.Range("B7").Value = "speridian.com"

# Test-data row (row 10): the short suffix code and the sample
# Fname/Uname values are refreshed to the new "aditya" persona.
$ws.Range("B10").Value = "n"
$ws.Range("J10").Value = "aditya"
$ws.Range("M10").Value = "aditya1"

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("M10").Select()
